$p = $ppt.ActivePresentation

# Slide 4 - "Example 1" TextBox -> "Solution"
$shp = $p.Slides.Item(4).Shapes.Item(4)
$shp.Left = 316.3782043457031
$shp.Width = 87.24378204345703
$shp.TextFrame.TextRange.Text = "Solution"

# Slide 7 - "Example 2" TextBox -> "Solution"
$shp = $p.Slides.Item(7).Shapes.Item(5)
$shp.Left = 464.5830993652344
$shp.Width = 87.24378204345703
$shp.TextFrame.TextRange.Text = "Solution"

# Slide 9 - "Example 3" TextBox -> "Solution"
$shp = $p.Slides.Item(9).Shapes.Item(4)
$shp.Left = 316.3782043457031
$shp.Width = 87.24378204345703
$shp.TextFrame.TextRange.Text = "Solution"

# Slide 11 - "Example 4" TextBox -> "Solution"
$shp = $p.Slides.Item(11).Shapes.Item(5)
$shp.Left = 451.003173828125
$shp.Width = 87.24378204345703
$shp.TextFrame.TextRange.Text = "Solution"
